$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = @"
Pipeline(steps=[('scaler',
                 ColumnTransformer(n_jobs=-1, remainder='passthrough',
                                   transformers=[('StandardScaler',
                                                  StandardScaler(),
                                                  ['AE_HR', 'AE_V',
                                                   'AbsOffAxis_HR',
                                                   'AbsOffAxis_V',
                                                   'AbsOnAxis_HR',
                                                   'AbsOnAxis_V', 'BallPath_HR',
                                                   'BallPath_V', 'CMT_HR',
                                                   'CMT_V', 'Corrective_HR',
                                                   'Corrective_V', 'Delta_AE',
                                                   'Delta_Fullpath', 'Delta_MT',
                                                   'Delta_OffAxis',
                                                   'Delta_OnAxis', 'D...
                                                   'Delta_RT', 'FullPath_HR',
                                                   'FullPath_V', 'MT_HR',
                                                   'MT_V', 'PeakV_HR',
                                                   'PeakV_V', 'RT_HR', 'RT_V',
                                                   'TMT_HR', 'TMT_V', 'VE_HR', ...])])),
                ('selector', 'passthrough'),
                ('model',
                 LGBMClassifier(boosting_type='dart', class_weight='balanced',
                                learning_rate=0.08362420395484067, max_bin=22,
                                max_depth=3, min_child_samples=16,
                                min_data_in_leaf=22, num_leaves=10,
                                random_state=42))])
"@
$ws.Range("C2").Value = @"
{'scaler': ColumnTransformer(n_jobs=-1, remainder='passthrough',
                  transformers=[('StandardScaler', StandardScaler(),
                                 ['AE_HR', 'AE_V', 'AbsOffAxis_HR',
                                  'AbsOffAxis_V', 'AbsOnAxis_HR', 'AbsOnAxis_V',
                                  'BallPath_HR', 'BallPath_V', 'CMT_HR',
                                  'CMT_V', 'Corrective_HR', 'Corrective_V',
                                  'Delta_AE', 'Delta_Fullpath', 'Delta_MT',
                                  'Delta_OffAxis', 'Delta_OnAxis', 'Delta_PV',
                                  'Delta_RT', 'FullPath_HR', 'FullPath_V',
                                  'MT_HR', 'MT_V', 'PeakV_HR', 'PeakV_V',
                                  'RT_HR', 'RT_V', 'TMT_HR', 'TMT_V', 'VE_HR', ...])]), 'model__boosting_type': 'dart', 'model__max_bin': 22, 'model__num_leaves': 10, 'model__learning_rate': 0.08362420395484067, 'model__max_depth': 3, 'model__min_data_in_leaf': 22, 'model__min_child_samples': 16, 'model__class_weight': 'balanced'}
"@
$ws.Range("D2").Value = 0.6178675297870344
$ws.Range("G2").Value = 0.8920213246588958
$ws.Range("H2").Value = 0.6666666666666666
$ws.Range("I2").Value = "[1 0 1 0 0 0 1 1 1 1 1 1 1 0 0 0 0 0 0 0 1 1 0 0]"
$ws.Range("J2").Value = "[0 0 1 0 0 1 0 1 1 0 1 1 1 0 0 1 1 1 0 0 1 1 0 1]"
$ws.Rows(2).AutoFit()

# Row 3
$ws.Range("B3").Value = @"
Pipeline(steps=[('scaler',
                 ColumnTransformer(n_jobs=-1, remainder='passthrough',
                                   transformers=[('StandardScaler',
                                                  StandardScaler(),
                                                  ['AE_HR', 'AE_V',
                                                   'AbsOffAxis_HR',
                                                   'AbsOffAxis_V',
                                                   'AbsOnAxis_HR',
                                                   'AbsOnAxis_V', 'BallPath_HR',
                                                   'BallPath_V', 'CMT_HR',
                                                   'CMT_V', 'Corrective_HR',
                                                   'Corrective_V', 'Delta_AE',
                                                   'Delta_Fullpath', 'Delta_MT',
                                                   'Delta_OffAxis',
                                                   'Delta_OnAxis', 'D...',
                                                   'Delta_RT', 'FullPath_HR',
                                                   'FullPath_V', 'MT_HR',
                                                   'MT_V', 'PeakV_HR',
                                                   'PeakV_V', 'RT_HR', 'RT_V',
                                                   'TMT_HR', 'TMT_V', 'VE_HR', ...])])),
                ('selector', 'passthrough'),
                ('model',
                 LGBMClassifier(boosting_type='dart', class_weight='balanced',
                                learning_rate=0.05784389350400665, max_bin=42,
                                max_depth=4, min_child_samples=6,
                                min_data_in_leaf=35, num_leaves=5,
                                random_state=42))])
"@
$ws.Range("C3").Value = @"
{'scaler': ColumnTransformer(n_jobs=-1, remainder='passthrough',
                  transformers=[('StandardScaler', StandardScaler(),
                                 ['AE_HR', 'AE_V', 'AbsOffAxis_HR',
                                  'AbsOffAxis_V', 'AbsOnAxis_HR', 'AbsOnAxis_V',
                                  'BallPath_HR', 'BallPath_V', 'CMT_HR',
                                  'CMT_V', 'Corrective_HR', 'Corrective_V',
                                  'Delta_AE', 'Delta_Fullpath', 'Delta_MT',
                                  'Delta_OffAxis', 'Delta_OnAxis', 'Delta_PV',
                                  'Delta_RT', 'FullPath_HR', 'FullPath_V',
                                  'MT_HR', 'MT_V', 'PeakV_HR', 'PeakV_V',
                                  'RT_HR', 'RT_V', 'TMT_HR', 'TMT_V', 'VE_HR', ...])]), 'model__boosting_type': 'dart', 'model__max_bin': 42, 'model__num_leaves': 5, 'model__learning_rate': 0.05784389350400665, 'model__max_depth': 4, 'model__min_data_in_leaf': 35, 'model__min_child_samples': 6, 'model__class_weight': 'balanced'}
"@
$ws.Range("D3").Value = 0.603634266340889
$ws.Range("G3").Value = 0.7420548476619618
$ws.Range("H3").Value = 0.542463768115942
$ws.Range("I3").Value = "[0 1 1 0 1 0 0 0 1 1 1 0 0 0 1 0 1 0 1 1 0 0 1 0]"
$ws.Range("J3").Value = "[0 1 1 1 0 0 0 0 1 1 1 1 1 1 0 1 0 0 1 0 0 0 0 1]"
$ws.Rows(3).AutoFit()

# Row 4
$ws.Range("B4").Value = @"
Pipeline(steps=[('scaler',
                 ColumnTransformer(n_jobs=-1, remainder='passthrough',
                                   transformers=[('StandardScaler',
                                                  StandardScaler(),
                                                  ['AE_HR', 'AE_V',
                                                   'AbsOffAxis_HR',
                                                   'AbsOffAxis_V',
                                                   'AbsOnAxis_HR',
                                                   'AbsOnAxis_V', 'BallPath_HR',
                                                   'BallPath_V', 'CMT_HR',
                                                   'CMT_V', 'Corrective_HR',
                                                   'Corrective_V', 'Delta_AE',
                                                   'Delta_Fullpath', 'Delta_MT',
                                                   'Delta_OffAxis',
                                                   'Delta_OnAxis', 'D...,
                                                   'Delta_RT', 'FullPath_HR',
                                                   'FullPath_V', 'MT_HR',
                                                   'MT_V', 'PeakV_HR',
                                                   'PeakV_V', 'RT_HR', 'RT_V',
                                                   'TMT_HR', 'TMT_V', 'VE_HR', ...])])),
                ('selector', 'passthrough'),
                ('model',
                 LGBMClassifier(boosting_type='dart', class_weight='balanced',
                                learning_rate=0.21911377955371011, max_bin=12,
                                max_depth=5, min_child_samples=11,
                                min_data_in_leaf=25, num_leaves=6,
                                random_state=42))])
"@
$ws.Range("C4").Value = @"
{'scaler': ColumnTransformer(n_jobs=-1, remainder='passthrough',
                  transformers=[('StandardScaler', StandardScaler(),
                                 ['AE_HR', 'AE_V', 'AbsOffAxis_HR',
                                  'AbsOffAxis_V', 'AbsOnAxis_HR', 'AbsOnAxis_V',
                                  'BallPath_HR', 'BallPath_V', 'CMT_HR',
                                  'CMT_V', 'Corrective_HR', 'Corrective_V',
                                  'Delta_AE', 'Delta_Fullpath', 'Delta_MT',
                                  'Delta_OffAxis', 'Delta_OnAxis', 'Delta_PV',
                                  'Delta_RT', 'FullPath_HR', 'FullPath_V',
                                  'MT_HR', 'MT_V', 'PeakV_HR', 'PeakV_V',
                                  'RT_HR', 'RT_V', 'TMT_HR', 'TMT_V', 'VE_HR', ...])]), 'model__boosting_type': 'dart', 'model__max_bin': 12, 'model__num_leaves': 6, 'model__learning_rate': 0.21911377955371011, 'model__max_depth': 5, 'model__min_data_in_leaf': 25, 'model__min_child_samples': 11, 'model__class_weight': 'balanced'}
"@
$ws.Range("D4").Value = 0.672122954444936
$ws.Range("G4").Value = 0.8925477897252091
$ws.Range("H4").Value = 0.4555072463768116
$ws.Range("I4").Value = "[0 1 1 0 0 1 0 0 0 1 0 1 0 1 0 1 0 0 1 0 0 1 1 1]"
$ws.Range("J4").Value = "[0 0 1 1 0 1 1 1 0 0 1 1 1 1 0 1 0 1 1 1 1 0 0 0]"
$ws.Rows(4).AutoFit()

# Row 5
$ws.Range("B5").Value = @"
Pipeline(steps=[('scaler',
                 ColumnTransformer(n_jobs=-1, remainder='passthrough',
                                   transformers=[('StandardScaler',
                                                  StandardScaler(),
                                                  ['AE_HR', 'AE_V',
                                                   'AbsOffAxis_HR',
                                                   'AbsOffAxis_V',
                                                   'AbsOnAxis_HR',
                                                   'AbsOnAxis_V', 'BallPath_HR',
                                                   'BallPath_V', 'CMT_HR',
                                                   'CMT_V', 'Corrective_HR',
                                                   'Corrective_V', 'Delta_AE',
                                                   'Delta_Fullpath', 'Delta_MT',
                                                   'Delta_OffAxis',
                                                   'Delta_OnAxis', 'D...,
                                                   'Delta_RT', 'FullPath_HR',
                                                   'FullPath_V', 'MT_HR',
                                                   'MT_V', 'PeakV_HR',
                                                   'PeakV_V', 'RT_HR', 'RT_V',
                                                   'TMT_HR', 'TMT_V', 'VE_HR', ...])])),
                ('selector', 'passthrough'),
                ('model',
                 LGBMClassifier(boosting_type='dart', class_weight='balanced',
                                learning_rate=0.02351732589932983, max_bin=22,
                                max_depth=2, min_child_samples=6,
                                min_data_in_leaf=22, num_leaves=10,
                                random_state=42))])
"@
$ws.Range("C5").Value = @"
{'scaler': ColumnTransformer(n_jobs=-1, remainder='passthrough',
                  transformers=[('StandardScaler', StandardScaler(),
                                 ['AE_HR', 'AE_V', 'AbsOffAxis_HR',
                                  'AbsOffAxis_V', 'AbsOnAxis_HR', 'AbsOnAxis_V',
                                  'BallPath_HR', 'BallPath_V', 'CMT_HR',
                                  'CMT_V', 'Corrective_HR', 'Corrective_V',
                                  'Delta_AE', 'Delta_Fullpath', 'Delta_MT',
                                  'Delta_OffAxis', 'Delta_OnAxis', 'Delta_PV',
                                  'Delta_RT', 'FullPath_HR', 'FullPath_V',
                                  'MT_HR', 'MT_V', 'PeakV_HR', 'PeakV_V',
                                  'RT_HR', 'RT_V', 'TMT_HR', 'TMT_V', 'VE_HR', ...])]), 'model__boosting_type': 'dart', 'model__max_bin': 22, 'model__num_leaves': 10, 'model__learning_rate': 0.02351732589932983, 'model__max_depth': 2, 'model__min_data_in_leaf': 22, 'model__min_child_samples': 6, 'model__class_weight': 'balanced'}
"@
$ws.Range("D5").Value = 0.5983189431378959
$ws.Range("G5").Value = 0.7956989247311828
$ws.Range("H5").Value = 0.4999999999999998
$ws.Range("I5").Value = "[0 0 0 1 1 1 1 1 0 0 1 0 1 0 0 0 1 0 1 1 0 0 1 0]"
$ws.Range("J5").Value = "[1 0 1 0 0 1 1 1 0 1 0 0 1 1 0 1 1 1 1 0 0 1 0 0]"
$ws.Rows(5).AutoFit()

# Row 6
$ws.Range("B6").Value = @"
Pipeline(steps=[('scaler',
                 ColumnTransformer(n_jobs=-1, remainder='passthrough',
                                   transformers=[('StandardScaler',
                                                  StandardScaler(),
                                                  ['AE_HR', 'AE_V',
                                                   'AbsOffAxis_HR',
                                                   'AbsOffAxis_V',
                                                   'AbsOnAxis_HR',
                                                   'AbsOnAxis_V', 'BallPath_HR',
                                                   'BallPath_V', 'CMT_HR',
                                                   'CMT_V', 'Corrective_HR',
                                                   'Corrective_V', 'Delta_AE',
                                                   'Delta_Fullpath', 'Delta_MT',
                                                   'Delta_OffAxis',
                                                   'Delta_OnAxis', 'D...',
                                                   'Delta_RT', 'FullPath_HR',
                                                   'FullPath_V', 'MT_HR',
                                                   'MT_V', 'PeakV_HR',
                                                   'PeakV_V', 'RT_HR', 'RT_V',
                                                   'TMT_HR', 'TMT_V', 'VE_HR', ...])])),
                ('selector', 'passthrough'),
                ('model',
                 LGBMClassifier(boosting_type='dart', class_weight='balanced',
                                learning_rate=0.48834529362316026, max_bin=32,
                                max_depth=6, min_child_samples=6,
                                min_data_in_leaf=22, num_leaves=2,
                                random_state=42))])
"@
$ws.Range("C6").Value = @"
{'scaler': ColumnTransformer(n_jobs=-1, remainder='passthrough',
                  transformers=[('StandardScaler', StandardScaler(),
                                 ['AE_HR', 'AE_V', 'AbsOffAxis_HR',
                                  'AbsOffAxis_V', 'AbsOnAxis_HR', 'AbsOnAxis_V',
                                  'BallPath_HR', 'BallPath_V', 'CMT_HR',
                                  'CMT_V', 'Corrective_HR', 'Corrective_V',
                                  'Delta_AE', 'Delta_Fullpath', 'Delta_MT',
                                  'Delta_OffAxis', 'Delta_OnAxis', 'Delta_PV',
                                  'Delta_RT', 'FullPath_HR', 'FullPath_V',
                                  'MT_HR', 'MT_V', 'PeakV_HR', 'PeakV_V',
                                  'RT_HR', 'RT_V', 'TMT_HR', 'TMT_V', 'VE_HR', ...])]), 'model__boosting_type': 'dart', 'model__max_bin': 32, 'model__num_leaves': 2, 'model__learning_rate': 0.48834529362316026, 'model__max_depth': 6, 'model__min_data_in_leaf': 22, 'model__min_child_samples': 6, 'model__class_weight': 'balanced'}
"@
$ws.Range("D6").Value = 0.6063996443872606
$ws.Range("G6").Value = 0.9569192298074519
$ws.Range("H6").Value = 0.3760869565217391
$ws.Range("I6").Value = "[1 0 1 0 1 1 0 0 0 1 0 1 0 1 1 1 0 1 0 1 0 0 0 0]"
$ws.Range("J6").Value = "[0 0 0 1 1 0 1 0 1 0 1 1 0 1 0 0 0 1 0 0 1 1 1 1]"
$ws.Rows(6).AutoFit()
